# Add a new header row at the top of Sheet1, labeling the existing two
# blocks of percentile columns (B:J and K:S) as "45_days" and "15_days"
# respectively. This pushes the old header/data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 1; everything that used to be row 1 (the old
# "stock" header row) and below shifts down to row 2, row 3, etc.
$ws.Range("A1").EntireRow.Insert()

# New row 1 labels: A1 stays empty, B1:J1 = "45_days", K1:S1 = "15_days".
$ws.Range("B1:J1").Value = "45_days"
$ws.Range("K1:S1").Value = "15_days"

Write-Output "Inserted header row with 45_days / 15_days labels"
